# Insert a new weekly price record above row 64 (shifts the existing
# rows 64..169 down to 65..170, matching the "Fruta / hortaliza, semanal"
# commit that adds one more week of data to the series).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 64 downward by inserting a blank row at position 64.
$ws.Rows.Item(64).Insert()

# Populate the newly inserted row 64 with the new weekly record. All
# fields other than the date (D) and volume (J) repeat the values of the
# record that used to sit in row 64 (now shifted to row 65).
$ws.Cells.Item(64, 1).Value = 7
$ws.Cells.Item(64, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(64, 3).Value = "Ñuble"
$ws.Cells.Item(64, 4).Value = 44580
$ws.Cells.Item(64, 5).Value = 16
$ws.Cells.Item(64, 6).Value = 100112017
$ws.Cells.Item(64, 7).Value = "Apio"
$ws.Cells.Item(64, 8).Value = "Americana (o)"
$ws.Cells.Item(64, 9).Value = "Primera"
$ws.Cells.Item(64, 10).Value = 60
$ws.Cells.Item(64, 11).Value = 8000
$ws.Cells.Item(64, 12).Value = 8500
$ws.Cells.Item(64, 13).Value = 8250
$ws.Cells.Item(64, 14).Value = "`$/docena de matas"
$ws.Cells.Item(64, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(64, 16).Value = 1375
$ws.Cells.Item(64, 17).Value = 6
$ws.Cells.Item(64, 18).Value = "Hortaliza"
